$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Supply" column (D) to "Y" and "Replacement" column (E)
# with the new supply-teacher-list names.
$ws.Range("D2").Value = "Y"
$ws.Range("E2").Value = "Dineth"

$ws.Range("D3").Value = "Y"
$ws.Range("E3").Value = "Marno"

$ws.Range("D4").Value = "Y"
$ws.Range("E4").Value = "Phillip"

$ws.Range("D5").Value = "Y"
$ws.Range("E5").Value = "Abdel"

# Widen column B slightly.
$ws.Columns.Item(2).ColumnWidth = 16

# Bold the header row (this also seeds the dxf used by the table's
# header-row style below).
$ws.Range("A1:E1").Font.Bold = $true

# Convert the range into a proper Excel Table ("Table2"), matching the
# "Supply Teacher List" table added to the workbook.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:E28"), $null, 1)
$tbl.Name = "Table2"
$tbl.TableStyle = "TableStyleLight16"

# Move the active selection.
[void]$ws.Range("P8").Select()
